# ---------------------------------------------------------------------------
# Rebuild the "URLs" sheet (was just a lone "Field" cell in A2) into a full
# VirusTotal-style results table, and normalise the "Hashes" sheet's
# Engine_detected (column D) values from numeric to text.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "URLs" worksheet
# ---------------------------------------------------------------------------
$urls = $wb.Worksheets.Item("URLs")

# Header row (row 1)
$urlHeaders = @("Field", "Type", "detected_url", "detected_urls_positives", "detected_urls_total", "detected_urls_scan_date", "Engine", "Engine_detected", "Engine_result")
for ($col = 1; $col -le $urlHeaders.Length; $col++) {
    $urls.Cells.Item(1, $col).Value = $urlHeaders[$col - 1]
}

# Data rows 2-7: one row per antivirus engine result for http://hecs.com
$engineRows = @(
    @("Comodo Valkyrie Verdict", "Comodo Valkyrie Verdict", "suspicious"),
    @("CRDF", "CRDF", "malicious"),
    @("Fortinet", "Fortinet", "malware"),
    @("AutoShun", "AutoShun", "malicious"),
    @("Webroot", "Webroot", "malicious"),
    @("Forcepoint ThreatSeeker", "Forcepoint ThreatSeeker", "suspicious")
)

$row = 2
foreach ($eng in $engineRows) {
    $urls.Cells.Item($row, 1).Value = "http://hecs.com"
    $urls.Cells.Item($row, 2).Value = "URL"
    $urls.Cells.Item($row, 3).Value = "http://hecs.com/"
    $urls.Cells.Item($row, 4).Value = 4
    $urls.Cells.Item($row, 5).Value = 86
    $urls.Cells.Item($row, 6).Value = "2022-05-31 10:11:18"
    $urls.Cells.Item($row, 7).Value = $eng[0]
    $urls.Cells.Item($row, 8).Value = $eng[1]
    $urls.Cells.Item($row, 9).Value = $eng[2]
    $row = $row + 1
}

# Styling: header row (A1:I1) and the "Field" key column (A1:A7) are bold,
# bordered, center/top aligned -- matching the style already used on A2's
# lone "Field" cell before this edit.
$urls.Range("A1:I1").Font.Bold = $true
$urls.Range("A1:I1").HorizontalAlignment = -4108
$urls.Range("A1:I1").VerticalAlignment = -4160
$urls.Range("A1:I1").Borders.LineStyle = 1

$urls.Range("A1:A7").Font.Bold = $true
$urls.Range("A1:A7").HorizontalAlignment = -4108
$urls.Range("A1:A7").VerticalAlignment = -4160
$urls.Range("A1:A7").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 2) "Hashes" worksheet -- Engine_detected (column D) numeric -> text
# ---------------------------------------------------------------------------
$hashes = $wb.Worksheets.Item("Hashes")

$hashes.Range("D2:D111").NumberFormat = "@"
for ($r = 2; $r -le 111; $r++) {
    $hashes.Cells.Item($r, 4).Value = "55"
}

$hashes.Range("D112:D173").NumberFormat = "@"
for ($r = 112; $r -le 173; $r++) {
    $hashes.Cells.Item($r, 4).Value = "62"
}
